$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement (Estado de Cuenta) rows 16-31 list "Periodo Mora" (col E)
# and "Valor Mora" (col F) entries. The old periods are being phased out and new
# ones added; the net effect on this sheet's data is that the period/value pairs
# for rows 16..31 end up in reverse order (oldest period moves to the bottom,
# newest period moves to the top), while every other column (doc type, doc
# number, worker name, salary, etc.) stays attached to its original row.

$firstRow = 16
$lastRow = 31
$periodCol = 5   # E: Periodo Mora
$valueCol = 6    # F: Valor Mora

$periods = @()
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, $periodCol).Text
    $values += $ws.Cells.Item($r, $valueCol).Value2
}

$count = $periods.Count
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $count - 1 - $i
    $ws.Cells.Item($targetRow, $periodCol).Value = $periods[$sourceIndex]
    $ws.Cells.Item($targetRow, $valueCol).Value = $values[$sourceIndex]
}
